# Update "want to go" counts (column F) for a handful of events that
# appear both on the "展览" sheet and on the combined "全部类型" sheet.
#
#   南宁·第一届ANE·DACG动漫嘉年华（取消）   1051 -> 1050
#   南宁·三月三漫次元国风动漫节             348  -> 352
#   南宁·2024三月三国潮动漫节（良牙春典）   2916 -> 2929
#   南宁·布谷鸟动漫展4th                    614  -> 615

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1050
$wsExpo.Range("F3").Value = 352
$wsExpo.Range("F4").Value = 2929
$wsExpo.Range("F6").Value = 615

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1050
$wsAll.Range("F5").Value = 352
$wsAll.Range("F6").Value = 2929
$wsAll.Range("F8").Value = 615
